$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NULL"
$ws.Range("B3").Value = "NULL"
$ws.Range("B4").Value = "NULL"
$ws.Range("B5").Value = "NULL"

$ws.Range("B5").Select()
